# Auto-generated edit script: applies numeric updates to H:N profit/price columns
# across several worksheets, per the target diff.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4833.1665
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 5249.75
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 5249.75
$ws.Range("M51").Value = -3516
$ws.Range("N51").Value = -6217.75
$ws.Range("H64").Value = 7642.143
$ws.Range("J64").Value = 7642.143
$ws.Range("L64").Value = 7642.143
$ws.Range("N64").Value = -8138.143
$ws.Range("H67").Value = 7642.143
$ws.Range("J67").Value = 7642.143
$ws.Range("L67").Value = 7642.143
$ws.Range("N67").Value = -9358.143
$ws.Range("H100").Value = 2684.9092
$ws.Range("I100").Value = 1764.2273
$ws.Range("K100").Value = 1764.2273
$ws.Range("M100").Value = -1223.2273
$ws.Range("H132").Value = 13730.862
$ws.Range("I132").Value = 14941.277
$ws.Range("J132").Value = 11750.182
$ws.Range("K132").Value = 44823.831
$ws.Range("L132").Value = 35250.546
$ws.Range("M132").Value = -42293.831
$ws.Range("N132").Value = -40310.546
$ws.Range("H138").Value = 2419.7666
$ws.Range("I138").Value = 3180.3684
$ws.Range("J138").Value = 2067.2927
$ws.Range("K138").Value = 9541.1052
$ws.Range("L138").Value = 6201.8781
$ws.Range("M138").Value = -4401.1052
$ws.Range("N138").Value = -16481.8781
$ws.Range("H141").Value = 4872.5713
$ws.Range("I141").Value = 4649.1816
$ws.Range("K141").Value = 13947.5448
$ws.Range("M141").Value = -8767.5448

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1671340.6
$ws.Range("I11").Value = 649.5
$ws.Range("J11").Value = 2506686.2
$ws.Range("K11").Value = 649.5
$ws.Range("L11").Value = 2506686.2
$ws.Range("M11").Value = -505.5
$ws.Range("N11").Value = -2506974.2
$ws.Range("H13").Value = 3333933.2
$ws.Range("I13").Value = 1300
$ws.Range("J13").Value = 5000250
$ws.Range("K13").Value = 1300
$ws.Range("L13").Value = 5000250
$ws.Range("M13").Value = -1156
$ws.Range("N13").Value = -5000538
$ws.Range("H74").Value = 13842.05
$ws.Range("I74").Value = 1587.16
$ws.Range("J74").Value = 34266.867
$ws.Range("K74").Value = 1587.16
$ws.Range("L74").Value = 34266.867
$ws.Range("M74").Value = -713.1600000000001
$ws.Range("N74").Value = -36014.867
$ws.Range("H77").Value = 13842.05
$ws.Range("I77").Value = 1587.16
$ws.Range("J77").Value = 34266.867
$ws.Range("K77").Value = 7935.8
$ws.Range("L77").Value = 171334.335
$ws.Range("M77").Value = -3567.8
$ws.Range("N77").Value = -180070.335
$ws.Range("H80").Value = 44999.5
$ws.Range("J80").Value = 44999.5
$ws.Range("L80").Value = 44999.5
$ws.Range("N80").Value = -46995.5
$ws.Range("H83").Value = 44999.5
$ws.Range("J83").Value = 44999.5
$ws.Range("L83").Value = 134998.5
$ws.Range("N83").Value = -144982.5
$ws.Range("H109").Value = 38251.332
$ws.Range("J109").Value = 38251.332
$ws.Range("L109").Value = 38251.332
$ws.Range("N109").Value = -41025.332
$ws.Range("H132").Value = 4421050
$ws.Range("I132").Value = 1995.1765
$ws.Range("J132").Value = 10199814
$ws.Range("K132").Value = 5985.529500000001
$ws.Range("L132").Value = 30599442
$ws.Range("M132").Value = -3455.529500000001
$ws.Range("N132").Value = -30604502

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H22").Value = 1267.5
$ws.Range("I22").Value = 1386.1111
$ws.Range("K22").Value = 1386.1111
$ws.Range("M22").Value = -1213.1111
$ws.Range("H134").Value = 100634.695
$ws.Range("I134").Value = 203959
$ws.Range("K134").Value = 611877
$ws.Range("M134").Value = -609342
$ws.Range("H135").Value = 85000
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 21906666
$ws.Range("I132").Value = 2319.6924
$ws.Range("J132").Value = 128690350
$ws.Range("K132").Value = 6959.0772
$ws.Range("L132").Value = 386071050
$ws.Range("M132").Value = -4429.0772
$ws.Range("N132").Value = -386076110
$ws.Range("H134").Value = 22228114
$ws.Range("I134").Value = 3843.724
$ws.Range("J134").Value = 62509604
$ws.Range("K134").Value = 11531.172
$ws.Range("L134").Value = 187528812
$ws.Range("M134").Value = -8996.172
$ws.Range("N134").Value = -187533882

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = 2100
$ws.Range("H139").Value = 47046
$ws.Range("I139").Value = 47046
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 141138
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -135998
$ws.Range("N139").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1711931.4
$ws.Range("I132").Value = 9592.571
$ws.Range("J132").Value = 5684055.5
$ws.Range("K132").Value = 28777.713
$ws.Range("L132").Value = 17052166.5
$ws.Range("M132").Value = -26247.713
$ws.Range("N132").Value = -17057226.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1181.75
$ws.Range("J10").Value = 1181.75
$ws.Range("L10").Value = 1181.75
$ws.Range("N10").Value = -1461.75
$ws.Range("H12").Value = 3840.6
$ws.Range("J12").Value = 3300.75
$ws.Range("L12").Value = 3300.75
$ws.Range("N12").Value = -3640.75
$ws.Range("H30").Value = 2429.8
$ws.Range("I30").Value = 1999.5
$ws.Range("J30").Value = 2716.6667
$ws.Range("K30").Value = 1999.5
$ws.Range("L30").Value = 2716.6667
$ws.Range("M30").Value = -1891.5
$ws.Range("N30").Value = -2932.6667
$ws.Range("H75").Value = 90173
$ws.Range("J75").Value = 90173
$ws.Range("L75").Value = 90173
$ws.Range("N75").Value = -92045
$ws.Range("H78").Value = 90173
$ws.Range("J78").Value = 90173
$ws.Range("L78").Value = 270519
$ws.Range("N78").Value = -279879

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 6733333
$ws.Range("J3").Value = 100000
$ws.Range("L3").Value = 100000
$ws.Range("N3").Value = -100228
$ws.Range("H10").Value = 7500
$ws.Range("I10").Value = 7500
$ws.Range("K10").Value = 7500
$ws.Range("M10").Value = -7331
$ws.Range("H75").Value = 270000
$ws.Range("J75").Value = 270000
$ws.Range("L75").Value = 270000
$ws.Range("N75").Value = -271872
$ws.Range("H78").Value = 270000
$ws.Range("J78").Value = 270000
$ws.Range("L78").Value = 810000
$ws.Range("N78").Value = -819360
$ws.Range("H132").Value = 478619.66
$ws.Range("I132").Value = 1899.35
$ws.Range("K132").Value = 5698.049999999999
$ws.Range("M132").Value = -3168.049999999999
$ws.Range("H136").Value = 527042.5600000001
$ws.Range("I136").Value = 2832.8462
$ws.Range("J136").Value = 1500574.9
$ws.Range("K136").Value = 8498.5386
$ws.Range("L136").Value = 4501724.699999999
$ws.Range("M136").Value = -5948.5386
$ws.Range("N136").Value = -4506824.699999999
